$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.130.33"
$ws.Range("E2").Value = "  +4.46%  "
$ws.Range("D3").Value = "2.246.83"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'244.82"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").Value = "'75.46"
$ws.Range("E7").Value = "  +7.85%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("E9").Value = "  +6.22%  "
$ws.Range("D10").Value = "'40.93"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").Value = "'0.0933"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "'6.98"
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "2.584.55"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "'14.62"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "2.244.50"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "'0.797"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "43.014.72"
$ws.Range("E18").Value = "  +4.68%  "
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("D20").Value = "'71.25"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'5.98"
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  +8.67%  "
$ws.Range("D23").Value = "'230.31"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  +16.22%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'10.89"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'3.46"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "'38.67"
$ws.Range("E29").Value = "  +28.35%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "'173.35"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").Value = "'20.30"
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("E33").Value = "  +5.31%  "
$ws.Range("D34").Value = "'5.31"
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  +7.04%  "
$ws.Range("E37").Value = "  +6.56%  "
$ws.Range("E38").Value = "  +19.53%  "
$ws.Range("D39").Value = "'13.10"
$ws.Range("E39").Value = "  +11.52%  "
$ws.Range("D41").Value = "'5.51"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("E42").Value = "  +8.35%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'59.76"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'105.42"
$ws.Range("E44").Value = "  +8.38%  "
$ws.Range("D45").Value = "'8.72"
$ws.Range("E45").Value = "  +5.44%  "
$ws.Range("D46").Value = "'0.488"
$ws.Range("E46").Value = "  +31.96%  "
$ws.Range("D47").Value = "'0.0993"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +9.57%  "
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").Value = "2.458.52"
$ws.Range("E51").Value = "  +3.40%  "

foreach ($addr in @("D5","D6","D7","D9","D10","D11","D12","D15","D17","D20","D21","D22","D23","D24","D26","D27","D28","D29","D31","D32","D34","D39","D41","D43","D44","D45","D46","D47","D48")) {
    $ws.Range($addr).Style = "Normal"
}
